# investing.com Mortgage Calculator functionality testing
# Adds a new "MortgageCalculator" worksheet (after "CalculatorCurrency"),
# populates it with the mortgage-calculator sample inputs/output, marks it
# as the active/selected sheet, and gives column A a best-fit width.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the current last sheet (CalculatorCurrency)
# so the tab order becomes: LogInValidCredentials, CalculatorCurrency, MortgageCalculator.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "MortgageCalculator"

# Populate the sample mortgage-calculator values (loan amount, interest rate,
# term in years, and the calculated monthly payment). Values are entered as
# text (leading apostrophe / quote-prefix) to match how the source sheet
# stores them.
$ws.Range("A3").Value = "'6"
$ws.Range("A2").Value = "'300000"
$ws.Range("A4").Value = "'20"
$ws.Range("A5").Value = "'$2,149.29"

# Force a text number format on the entered cells (keeps them quote-prefixed
# text rather than being reinterpreted as numbers/currency).
$ws.Range("A2:A5").NumberFormat = "@"

# Best-fit column A to its contents.
$ws.Columns("A:A").AutoFit() | Out-Null

# Select C6 and make this the active sheet/tab, matching the authored state.
$ws.Range("C6").Select() | Out-Null
